# Auto-generated edit script: refresh profit/cost computed columns (H-N)
# across multiple sheets, as produced by the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3650.5334
$ws.Range("J19").Value = 4251
$ws.Range("L19").Value = 4251
$ws.Range("N19").Value = -4601
$ws.Range("H40").Value = 2284.2307
$ws.Range("I40").Value = 1833.3334
$ws.Range("J40").Value = 2419.5
$ws.Range("K40").Value = 1833.3334
$ws.Range("L40").Value = 2419.5
$ws.Range("M40").Value = -1658.3334
$ws.Range("N40").Value = -2769.5
$ws.Range("H100").Value = 1405.24
$ws.Range("I100").Value = 1199.2632
$ws.Range("J100").Value = 2057.5
$ws.Range("K100").Value = 1199.2632
$ws.Range("L100").Value = 2057.5
$ws.Range("M100").Value = -658.2632000000001
$ws.Range("N100").Value = -3139.5
$ws.Range("H103").Value = 399
$ws.Range("I103").Value = 399.8
$ws.Range("J103").Value = 395
$ws.Range("K103").Value = 1199.4
$ws.Range("L103").Value = 1185
$ws.Range("M103").Value = -613.4000000000001
$ws.Range("N103").Value = -2357
$ws.Range("H129").Value = 2039.2
$ws.Range("J129").Value = 2039.2
$ws.Range("L129").Value = 6117.6
$ws.Range("N129").Value = -16117.6
$ws.Range("H131").Value = 8809.929
$ws.Range("I131").Value = 2977
$ws.Range("K131").Value = 8931
$ws.Range("M131").Value = -3891
$ws.Range("H132").Value = 2267.8953
$ws.Range("I132").Value = 2085.6184
$ws.Range("K132").Value = 6256.8552
$ws.Range("M132").Value = -3726.8552
$ws.Range("H137").Value = 27779844
$ws.Range("I137").Value = 62501576
$ws.Range("J137").Value = 2460.1
$ws.Range("K137").Value = 187504728
$ws.Range("L137").Value = 7380.299999999999
$ws.Range("M137").Value = -187502178
$ws.Range("N137").Value = -12480.3
$ws.Range("H141").Value = 1015.12
$ws.Range("I141").Value = 1026.1666
$ws.Range("J141").Value = 750
$ws.Range("K141").Value = 3078.4998
$ws.Range("L141").Value = 2250
$ws.Range("M141").Value = 2101.5002
$ws.Range("N141").Value = -12610

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1372.9333
$ws.Range("I61").Value = 1151.4546
$ws.Range("J61").Value = 1982
$ws.Range("K61").Value = 1151.4546
$ws.Range("L61").Value = 1982
$ws.Range("M61").Value = -939.4546
$ws.Range("N61").Value = -2406
$ws.Range("H74").Value = 5676.0586
$ws.Range("I74").Value = 6178.5713
$ws.Range("J74").Value = 3331
$ws.Range("K74").Value = 6178.5713
$ws.Range("L74").Value = 3331
$ws.Range("M74").Value = -5304.5713
$ws.Range("N74").Value = -5079
$ws.Range("H77").Value = 5676.0586
$ws.Range("I77").Value = 6178.5713
$ws.Range("J77").Value = 3331
$ws.Range("K77").Value = 30892.8565
$ws.Range("L77").Value = 16655
$ws.Range("M77").Value = -26524.8565
$ws.Range("N77").Value = -25391
$ws.Range("H97").Value = 647.1071
$ws.Range("I97").Value = 656.0417
$ws.Range("J97").Value = 593.5
$ws.Range("K97").Value = 656.0417
$ws.Range("L97").Value = 593.5
$ws.Range("M97").Value = -160.0417
$ws.Range("N97").Value = -1585.5
$ws.Range("H102").Value = 38457.777
$ws.Range("I102").Value = 29057.666
$ws.Range("K102").Value = 29057.666
$ws.Range("M102").Value = -27435.666
$ws.Range("H122").Value = 3036.25
$ws.Range("I122").Value = 1707.7333
$ws.Range("K122").Value = 5123.199900000001
$ws.Range("M122").Value = -2673.199900000001
$ws.Range("H125").Value = 50250
$ws.Range("J125").Value = 50250
$ws.Range("L125").Value = 50250
$ws.Range("N125").Value = -60090
$ws.Range("H132").Value = 1285.619
$ws.Range("I132").Value = 1245.4615
$ws.Range("J132").Value = 1350.875
$ws.Range("K132").Value = 3736.3845
$ws.Range("L132").Value = 4052.625
$ws.Range("M132").Value = -1206.3845
$ws.Range("N132").Value = -9112.625
$ws.Range("H136").Value = 1372.9333
$ws.Range("I136").Value = 1151.4546
$ws.Range("J136").Value = 1982
$ws.Range("K136").Value = 3454.3638
$ws.Range("L136").Value = 5946
$ws.Range("M136").Value = -904.3638000000001
$ws.Range("N136").Value = -11046

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7851
$ws.Range("I20").Value = 11199.25
$ws.Range("K20").Value = 11199.25
$ws.Range("M20").Value = -10952.25
$ws.Range("H22").Value = 1410.3334
$ws.Range("I22").Value = 115.5
$ws.Range("J22").Value = 4000
$ws.Range("K22").Value = 115.5
$ws.Range("L22").Value = 4000
$ws.Range("M22").Value = 57.5
$ws.Range("N22").Value = -4346
$ws.Range("H134").Value = 1023.44446
$ws.Range("I134").Value = 999.625
$ws.Range("K134").Value = 2998.875
$ws.Range("M134").Value = -463.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()
$ws.Range("H22").Value = 668.125
$ws.Range("I22").Value = 587.95
$ws.Range("J22").Value = 1069
$ws.Range("K22").Value = 587.95
$ws.Range("L22").Value = 1069
$ws.Range("M22").Value = -237.95
$ws.Range("N22").Value = -1769
$ws.Range("H58").Value = 1411.3636
$ws.Range("I58").Value = 1365.2667
$ws.Range("J58").Value = 1510.1428
$ws.Range("K58").Value = 1365.2667
$ws.Range("L58").Value = 1510.1428
$ws.Range("M58").Value = -1162.2667
$ws.Range("N58").Value = -1916.1428
$ws.Range("H99").Value = 2381.1765
$ws.Range("I99").Value = 2199
$ws.Range("K99").Value = 2199
$ws.Range("M99").Value = -701
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
$ws.Range("H126").Value = 2381.1765
$ws.Range("I126").Value = 2199
$ws.Range("K126").Value = 6597
$ws.Range("M126").Value = -4127
$ws.Range("H132").Value = 3439.4482
$ws.Range("I132").Value = 3023.1428
$ws.Range("K132").Value = 9069.4284
$ws.Range("M132").Value = -6539.428400000001
$ws.Range("H134").Value = 1655.1666
$ws.Range("I134").Value = 1543.7667
$ws.Range("J134").Value = 2212.1667
$ws.Range("K134").Value = 4631.300099999999
$ws.Range("L134").Value = 6636.500100000001
$ws.Range("M134").Value = -2096.300099999999
$ws.Range("N134").Value = -11706.5001
$ws.Range("H136").Value = 1411.3636
$ws.Range("I136").Value = 1365.2667
$ws.Range("J136").Value = 1510.1428
$ws.Range("K136").Value = 4095.800099999999
$ws.Range("L136").Value = 4530.428400000001
$ws.Range("M136").Value = -1545.800099999999
$ws.Range("N136").Value = -9630.4284
$ws.Range("H140").Value = 75427
$ws.Range("I140").Value = 20709
$ws.Range("J140").Value = 93666.336
$ws.Range("K140").Value = 20709
$ws.Range("L140").Value = 93666.336
$ws.Range("M140").Value = -15529
$ws.Range("N140").Value = -104026.336

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8953.875
$ws.Range("I70").Value = 9185.615
$ws.Range("K70").Value = 9185.615
$ws.Range("M70").Value = -8915.615
$ws.Range("H73").Value = 8953.875
$ws.Range("I73").Value = 9185.615
$ws.Range("K73").Value = 9185.615
$ws.Range("M73").Value = -8249.615
$ws.Range("H80").Value = 4045.762
$ws.Range("I80").Value = 3828.8667
$ws.Range("K80").Value = 3828.8667
$ws.Range("M80").Value = -2830.8667
$ws.Range("H83").Value = 4045.762
$ws.Range("I83").Value = 3828.8667
$ws.Range("K83").Value = 19144.3335
$ws.Range("M83").Value = -14152.3335
$ws.Range("H122").Value = 1468.4814
$ws.Range("I122").Value = 1245.409
$ws.Range("K122").Value = 3736.227
$ws.Range("M122").Value = -1286.227
$ws.Range("H132").Value = 7212.4375
$ws.Range("I132").Value = 6268.8887
$ws.Range("K132").Value = 18806.6661
$ws.Range("M132").Value = -16276.6661

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 787.5
$ws.Range("J22").Value = 1523.2
$ws.Range("L22").Value = 1523.2
$ws.Range("N22").Value = -2113.2
$ws.Range("H27").Value = 787.5
$ws.Range("J27").Value = 1523.2
$ws.Range("L27").Value = 1523.2
$ws.Range("N27").Value = -1737.2
$ws.Range("H46").Value = 2557.1875
$ws.Range("I46").Value = 1530.1875
$ws.Range("J46").Value = 3584.1875
$ws.Range("K46").Value = 1530.1875
$ws.Range("L46").Value = 3584.1875
$ws.Range("M46").Value = -1342.1875
$ws.Range("N46").Value = -3960.1875
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 30113
$ws.Range("J63").Value = 30000
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31248
$ws.Range("H66").Value = 30113
$ws.Range("J66").Value = 30000
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -96240
$ws.Range("H132").Value = 3252.0667
$ws.Range("I132").Value = 3141.1714
$ws.Range("K132").Value = 9423.514200000001
$ws.Range("M132").Value = -6893.514200000001
$ws.Range("H136").Value = 2709.077
$ws.Range("I136").Value = 1092.8
$ws.Range("J136").Value = 4913.091
$ws.Range("K136").Value = 3278.4
$ws.Range("L136").Value = 14739.273
$ws.Range("M136").Value = -728.3999999999996
$ws.Range("N136").Value = -19839.273
